# Fruta / hortaliza, semanal
# A new weekly price record (Acelga, Femacal de La Calera) is inserted at
# row 181, pushing the existing rows 181-195 down to 182-196.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 181, shifting the rest down.
$ws.Rows.Item(181).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(181, 1).Value = 3
$ws.Cells.Item(181, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(181, 3).Value = "Coquimbo"
$ws.Cells.Item(181, 4).Value = 44461
$ws.Cells.Item(181, 5).Value = 5
$ws.Cells.Item(181, 6).Value = 100112009
$ws.Cells.Item(181, 7).Value = "Acelga"
$ws.Cells.Item(181, 8).Value = "Sin especificar"
$ws.Cells.Item(181, 9).Value = "Primera"
$ws.Cells.Item(181, 10).Value = 260
$ws.Cells.Item(181, 11).Value = 2000
$ws.Cells.Item(181, 12).Value = 2200
$ws.Cells.Item(181, 13).Value = 2092
$ws.Cells.Item(181, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(181, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(181, 16).Value = 349
$ws.Cells.Item(181, 17).Value = 6
$ws.Cells.Item(181, 18).Value = "Hortaliza"
